# webapply postman collections and scenarios file
#
# Adds the new "SME Proposal" / "Search SME/Retail/Corporate Proposal"
# scenario rows (rows 9-13, plus the two new request/response cells on
# row 8) to the WiremockScenarios sheet, and refreshes the sheet view
# (selection/zoom) and a few column widths to match the post-edit
# state produced in Excel.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Row 8: add the request/response mapping file names for BackofficeLogin ---
$ws.Range("F8").Value = "COSME0007BackOfficeLoginRequest"
$ws.Range("G8").Value = "COSME0007BackOfficeLoginResponse"

# --- Row 9: SME Proposal ---
$ws.Range("C9").Value = "SMEProposal"
$ws.Range("D9").Value = "SMEProposal"
$ws.Range("E9").Value = "prospectId= COSME008"
$ws.Range("F9").Value = "COSME0008SMEProspectRequest"
$ws.Range("G9").Value = "COSME0008SMEProspectResponse"

# --- Row 10: Update SME Proposal ---
$ws.Range("C10").Value = "Update SME Proposal"
$ws.Range("D10").Value = "UpdateSMEProposal"
$ws.Range("E10").Value = "prospectId= COSME008"
$ws.Range("F10").Value = "COSME0009SearchSMEProspectRequest"
$ws.Range("G10").Value = "COSME0008SMEProspectResponse"

# --- Row 11: Search SME Proposal ---
$ws.Range("C11").Value = "Search SME Proposal"
$ws.Range("D11").Value = "SearchSMEProposal"
$ws.Range("E11").Value = "fname = COSME0010"
$ws.Range("F11").Value = "COSME0010SearchSMEProspectRequest"
$ws.Range("G11").Value = "COSME0010SearchSMEProspectResponse"

# --- Row 12: Search Retail Proposal ---
$ws.Range("C12").Value = "Search Retail Proposal"
$ws.Range("D12").Value = "SearchSMEProposal"
$ws.Range("E12").Value = "fname = COSME0010"
$ws.Range("F12").Value = "COSME0011SearchRetailProspectRequest"
$ws.Range("G12").Value = "COSME0010SearchSMEProspectResponse"

# --- Row 13: Search Corporate Proposal ---
$ws.Range("C13").Value = "Search Corporate Proposal"
$ws.Range("D13").Value = "SearchCorporateProposal"
$ws.Range("E13").Value = "fname = COSME0010"
$ws.Range("F13").Value = "COSME0012SearchCorporateProspectRequest"
$ws.Range("G13").Value = "COSME0010SearchSMEProspectResponse"

# --- Column width refresh (columns widened/narrowed to fit the new data) ---
$ws.Columns.Item(2).ColumnWidth = 10.333333333333332
$ws.Columns.Item(3).ColumnWidth = 52.33333333333333
$ws.Columns.Item(4).ColumnWidth = 21.166666666666668
$ws.Columns.Item(7).ColumnWidth = 34.666666666666664

# --- View state: zoom back to 100% and move the selection to C10 ---
$win = $excel.ActiveWindow
$win.Zoom = 100
$ws.Range("C10").Select()
